# "Terminado zip Entrega 2"
# Marks Fase 2 (Entregables + Rubrica) as done, rewrites the "Metodo" sheet
# progress tracker (bold phase headers, closes blank separator rows, adds
# "x" marks + Laberintos/guia/tesoros notes), switches the active sheet to
# "Metodo" with a new selection, and tweaks the other sheets' selections.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Entregables
$ws2 = $wb.Worksheets.Item(2)   # Rubrica
$ws3 = $wb.Worksheets.Item(3)   # Metodo

# --- Entregables: mark "Fase 2" of the "Informe metodo de la ingenieria" row done
$ws1.Range('C13').Value = 'x'

# --- Rubrica: mark "Fase 3. Busqueda de soluciones creativas" row done
$ws2.Range('C11').Value = 'x'

# --- Metodo: rewrite the tracker, closing up the blank spacer rows and
# marking each already-covered item with an "x"
$ws3.Range('A1').Value = 'Fase 1'

$ws3.Range('A5').Value = 'Fase 2'
$ws3.Range('C5').Value = 'x'
$ws3.Range('A6').Value = 'Contexto'
$ws3.Range('B6').Value = 'Laberintos'
$ws3.Range('C6').Value = 'x'
$ws3.Range('A7').Value = 'Teoría'
$ws3.Range('A8').Value = 'Grafos'
$ws3.Range('C8').Value = 'x'
$ws3.Range('A9').Value = 'BFS Y DFS'
$ws3.Range('A10').Value = 'Camino min'
$ws3.Range('A11').Value = 'Árboles de recubrimiento'

$ws3.Range('A12').Value = 'Fase 3'
$ws3.Range('A13').Value = 'Método de generación de ideas'
$ws3.Range('B13').Value = 'Para la guía'
$ws3.Range('C13').Value = 'x'
$ws3.Range('A14').Value = 'Al menos 7 ideas'
$ws3.Range('B14').Value = 'Para los tesoros'

$ws3.Range('A15').Value = 'Fase 4'
$ws3.Range('A16').Value = 'Descarte de ideas'
$ws3.Range('B16').Value = 'Para la guía'
$ws3.Range('C16').Value = 'x'
$ws3.Range('A17').Value = 'Diseño preliminar otras ideas'
$ws3.Range('B17').Value = 'Para los tesoros'

$ws3.Range('A18').Value = 'Fase 5'
$ws3.Range('A19').Value = 'Criterios'
$ws3.Range('B19').Value = 'Para la guía'
$ws3.Range('C19').Value = 'x'
$ws3.Range('A20').Value = 'Asignar valores '
$ws3.Range('B20').Value = 'Para los tesoros'

# cells left behind by the row shuffle that must now be empty
$ws3.Range('B7').Clear()
$ws3.Range('B15').Clear()
$ws3.Range('C15').Clear()
$ws3.Range('A22').Clear()
$ws3.Range('A23').Clear()
$ws3.Range('A24').Clear()

# bold the "Fase N" section headers in column A
$ws3.Range('A1').Font.Bold = $true
$ws3.Range('A5').Font.Bold = $true
$ws3.Range('A12').Font.Bold = $true
$ws3.Range('A15').Font.Bold = $true
$ws3.Range('A18').Font.Bold = $true

# page orientation for the Metodo sheet
$ws3.PageSetup.Orientation = 1

# selections on each sheet (Entregables and Rubrica are no longer the active tab)
$null = $ws1.Range('D10').Select()
$null = $ws2.Range('C12').Select()
$null = $ws3.Range('A3').Select()
